# Updates the "Price" (column D) and "Volume(1h)" (column E) figures on the
# cryptocurrency listing sheet, mirroring the periodic "Updated symbol list"
# GitHub Actions refresh of this workbook.
#
# Every value in these two columns is stored as literal text (e.g. "312.51",
# "-0.72%") rather than as a number/percentage, so each cell's number format
# is forced to Text ("@") before the new value is written. This stops Excel
# from re-interpreting strings such as "312.51" or "-0.72%" as a numeric
# value / percentage and silently changing the stored cell type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row new values for column D (Price) and column E (Volume(1h)).
# A value of $null means that column is left unchanged for that row.
$updates = @(
    @{ Row = 2;  D = '312.51';         E = '-0.72%' },
    @{ Row = 3;  D = '37.61';          E = '-4.18%' },
    @{ Row = 4;  D = '5.072';          E = '-1.42%' },
    @{ Row = 5;  D = '0.07755';        E = '-5.05%' },
    @{ Row = 6;  D = '4.358';          E = '-0.54%' },
    @{ Row = 7;  D = '1.890';          E = '-4.93%' },
    @{ Row = 8;  D = '8.197';          E = '-1.58%' },
    @{ Row = 9;  D = '2.964';          E = '-5.18%' },
    @{ Row = 10; D = '0.9170';         E = '-2.04%' },
    @{ Row = 11; D = '0.1227';         E = '-5.12%' },
    @{ Row = 12; D = '0.1900';         E = '-3.70%' },
    @{ Row = 13; D = '0.08886';        E = '-2.56%' },
    @{ Row = 14; D = '0.03395';        E = '-3.69%' },
    @{ Row = 15; D = '0.09703';        E = '-0.32%' },
    @{ Row = 16; D = '0.001378';       E = '-2.68%' },
    @{ Row = 17; D = '0.005724';       E = '-5.94%' },
    @{ Row = 18; D = '3.540';          E = '-2.40%' },
    @{ Row = 19; D = $null;            E = '-1.75%' },
    @{ Row = 20; D = '0.1287';         E = '-1.60%' },
    @{ Row = 21; D = '5.042';          E = '-0.55%' },
    @{ Row = 22; D = '0.2594';         E = '4.24%' },
    @{ Row = 23; D = '0.02107';        E = '5,595.88%' },
    @{ Row = 24; D = '0.04387';        E = '0.56%' },
    @{ Row = 25; D = '0.001215';       E = '-2.22%' },
    @{ Row = 26; D = '0.004235';       E = '-11.01%' },
    @{ Row = 27; D = '0.0001351';      E = '-65.28%' },
    @{ Row = 39; D = '0.02125';        E = '-5.20%' },
    @{ Row = 40; D = '0.04961';        E = '-4.75%' },
    @{ Row = 41; D = '0.007785';       E = '0.76%' },
    @{ Row = 42; D = '0.009964';       E = '-3.08%' },
    @{ Row = 43; D = $null;            E = '-4.05%' },
    @{ Row = 44; D = '0.001994';       E = '-5.10%' },
    @{ Row = 45; D = '0.009657';       E = '8.99%' },
    @{ Row = 46; D = '0.00006596';     E = '-3.27%' },
    @{ Row = 47; D = '0.00000000751';  E = '0.08%' },
    @{ Row = 48; D = '0.003046';       E = '1.27%' },
    @{ Row = 50; D = '0.00002102';     E = '0.08%' },
    @{ Row = 51; D = '0.0002002';      E = '0.08%' }
)

foreach ($update in $updates) {
    if ($null -ne $update.D) {
        $cell = $ws.Range("D" + $update.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $update.D
    }
    if ($null -ne $update.E) {
        $cell = $ws.Range("E" + $update.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $update.E
    }
}
